$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 42 and 43 need their identifying/text columns (A-D) swapped, and the
# Station18 (J) value moves from row 43 to row 42.

# Set row 42 to the values that used to be on row 43 (A-D), plus add J42 = 0
$ws.Range("A42").Value = "c0a3f3ed23f04247d92740a9502f8b57"
$ws.Range("B42").Value = "unassigned"
$ws.Range("C42").Value = "unassigned"
$ws.Range("D42").Value = "unassigned"
$ws.Range("J42").Value = 0

# Set row 43 to the values that used to be on row 42 (A-D), and clear J43
$ws.Range("A43").Value = "307c55294ffe3b8aa46fce358d55590e"
$ws.Range("B43").Value = "Homo sapiens"
$ws.Range("C43").Value = "Human"
$ws.Range("D43").Value = "Human"
$ws.Range("J43").ClearContents()
